# Bid progress migration & agregation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: task reassigned from Jakub Ivan Vanko to Tomáš Adam; description updated
$ws.Range("C6").Value = "Tomáš Adam"
$ws.Range("E6").Value = "Auction progress visualizaion & filtering by columns and participats agregation"

# Row 8: description gains a trailing space (author re-typed the cell)
$ws.Range("E8").Value = "Dashboard - overview and charts "

# Row 9: task removed entirely (Karin Jana Szilárdy / UI - Linear regression / Price prediction - Linear regression)
$ws.Range("C9").Value = ""
$ws.Range("D9").Value = ""
$ws.Range("E9").Value = ""

# Refresh selection / zoom to match author's final view state
$ws.Range("E8").Select() | Out-Null
$excel.ActiveWindow.Zoom = 114
